$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.893.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.597.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.77%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "

$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("E11").Value = "  +5.09%  "

$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.055.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.842.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.590.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  -2.00%  "

$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.41%  "

$ws.Range("E21").Value = "  -1.84%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0754"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.876"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.73%  "

$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "282.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.599"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0957"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.930.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "

$ws.Range("E50").Value = "  -2.91%  "

$ws.Range("E51").Value = "  -2.27%  "
